$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = 2..29

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$c]
    }
}

function Swap-Rows($rowA, $rowB) {
    $a = Get-RowValues $rowA
    $b = Get-RowValues $rowB
    Set-RowValues $rowA $b
    Set-RowValues $rowB $a
}

# Simple pairwise swaps of match records (B:AC), ids (col A) stay put.
Swap-Rows 27 28
Swap-Rows 47 48
Swap-Rows 104 105
Swap-Rows 148 149

# 3-way rotation among rows 142/143/144:
#   new142 <- old144, new143 <- old142, new144 <- old143
$v142 = Get-RowValues 142
$v143 = Get-RowValues 143
$v144 = Get-RowValues 144
Set-RowValues 142 $v144
Set-RowValues 143 $v142
Set-RowValues 144 $v143

# Append a brand-new match row (row 209), copying formatting from the last
# existing data row (208) so styles (A bold/bordered, E date format) match.
$ws.Range("A208:AC208").Copy() | Out-Null
$ws.Range("A209:AC209").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(209, 1).Value2 = 207
$ws.Cells.Item(209, 2).Value2 = 8010641
$ws.Cells.Item(209, 3).Value2 = "Bolivia Primera División"
$ws.Cells.Item(209, 4).Value2 = "Bolivia Apertura"
$ws.Cells.Item(209, 5).Value2 = 45383.875
$ws.Cells.Item(209, 6).Value2 = "Jorge Wilstermann"
$ws.Cells.Item(209, 7).Value2 = "San Jose de Oruro"
$ws.Cells.Item(209, 8).Value2 = 3
$ws.Cells.Item(209, 9).Value2 = 4
$ws.Cells.Item(209, 10).Value2 = "A"
$ws.Cells.Item(209, 11).Value2 = 1.7
$ws.Cells.Item(209, 12).Value2 = 3.4
$ws.Cells.Item(209, 13).Value2 = 4.333
$ws.Cells.Item(209, 14).Value2 = 1.5
$ws.Cells.Item(209, 15).Value2 = 4.2
$ws.Cells.Item(209, 16).Value2 = 7
$ws.Cells.Item(209, 17).Value2 = -1.25
$ws.Cells.Item(209, 18).Value2 = 1.975
$ws.Cells.Item(209, 19).Value2 = 1.825
$ws.Cells.Item(209, 20).Value2 = 2.75
$ws.Cells.Item(209, 21).Value2 = 1.95
$ws.Cells.Item(209, 22).Value2 = 1.85
$ws.Cells.Item(209, 23).Value2 = -1
$ws.Cells.Item(209, 24).Value2 = -1
$ws.Cells.Item(209, 25).Value2 = 6
$ws.Cells.Item(209, 26).Value2 = -1
$ws.Cells.Item(209, 27).Value2 = 0.825
$ws.Cells.Item(209, 28).Value2 = 0.95
$ws.Cells.Item(209, 29).Value2 = -1
